# Adds new sprint-planning rows (S20 / G01 - "Holdings trade sizing & risk controls")
# to the bottom of the tasks table on Sheet1, rows 167-170, columns A-H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Indian Rupee sign (U+20B9) built from its code point to avoid any
# source-file encoding ambiguity.
$rupee = [char]0x20B9

$rows = @(
    @{
        A = "S20"
        B = "G01"
        C = "Holdings trade sizing & risk controls"
        D = "S20_G01_TB001"
        E = "Expose a lightweight portfolio summary (live equity and per-symbol position values) to support % of portfolio and risk-based sizing calculations."
        F = "Reuses existing holdings/positions data; initial implementation can be an internal helper or small API endpoint without full reporting UI."
        G = "planned"
        H = "Provides a single source of truth for total portfolio value and per-position notional used by sizing modes."
    },
    @{
        A = "S20"
        B = "G01"
        C = "Holdings trade sizing & risk controls"
        D = "S20_G01_TF001"
        E = "Extend Buy/Sell dialog with unified sizing modes (Qty, Amount, % of position, % of portfolio) and consistent auto-calculation between them."
        F = "Builds on the current Qty/Amount/% of position behaviour, adding a % of portfolio mode that derives target notional from portfolio value."
        G = "planned"
        H = "Keeps day-to-day trading simple while enabling portfolio-level rebalancing directly from the holdings dialog."
    },
    @{
        A = "S20"
        B = "G01"
        C = "Holdings trade sizing & risk controls"
        D = "S20_G01_TB002"
        E = "Add backend helper to compute risk-based position size from entry price, stop level, and risk budget ($rupee or % of portfolio)."
        F = "Implemented as a pure function/service that can be reused later by strategies or analytics; no execution decisions are automated in this phase."
        G = "planned"
        H = "Encodes the core risk per share and max-loss sizing formulas in one validated place."
    },
    @{
        A = "S20"
        B = "G01"
        C = "Holdings trade sizing & risk controls"
        D = "S20_G01_TF002"
        E = "Introduce a Risk sizing mode in the Buy/Sell dialog that lets the user specify risk budget and stop price and shows derived qty/amount and expected max loss."
        F = "Risk mode is optional and advanced; it calls the backend helper for calculations but still requires the user to confirm and submit orders manually."
        G = "planned"
        H = "Makes it easier to keep per-trade downside consistent without changing the existing order-routing pipeline."
    }
)

$startRow = 167
$colOrder = @("A", "B", "C", "D", "E", "F", "G", "H")

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    foreach ($col in $colOrder) {
        $cell = $ws.Range("$col$r")
        $cell.Value2 = $rowData[$col]
        # New rows are appended without the table's usual wrap-text style,
        # matching the plain/default formatting used in the source edit.
        $cell.Style = "Normal"
    }
}
